# Updates 'F' column (想去人数 / interested-count) values on sheets
# 展览, 演出, and 全部类型 to the freshly scraped snapshot figures.
# Mirrors commit: 'Update gh-pages to output generated at 456a3b4'.

$wb = $excel.ActiveWorkbook

# Sheet 展览 (Exhibition) - update 想去人数 (interested count) column F
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 3013
$ws.Range("F3").Value = 3013
$ws.Range("F4").Value = 6434
$ws.Range("F5").Value = 2556
$ws.Range("F6").Value = 655
$ws.Range("F7").Value = 88
$ws.Range("F8").Value = 3147
$ws.Range("F9").Value = 370
$ws.Range("F11").Value = 7714
$ws.Range("F12").Value = 386
$ws.Range("F15").Value = 10
$ws.Range("F16").Value = 6
$ws.Range("F17").Value = 262
$ws.Range("F19").Value = 21
$ws.Range("F21").Value = 9487
$ws.Range("F23").Value = 271
$ws.Range("F28").Value = 132
$ws.Range("F29").Value = 36
$ws.Range("F33").Value = 2627
$ws.Range("F35").Value = 2054
$ws.Range("F36").Value = 18
$ws.Range("F37").Value = 1492
$ws.Range("F38").Value = 803
$ws.Range("F39").Value = 3974
$ws.Range("F40").Value = 221
$ws.Range("F41").Value = 617
$ws.Range("F42").Value = 1200
$ws.Range("F43").Value = 115
$ws.Range("F44").Value = 258
$ws.Range("F45").Value = 76
$ws.Range("F49").Value = 68
$ws.Range("F50").Value = 16

# Sheet 演出 (Performance) - update 想去人数 column F
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 1
$ws.Range("F8").Value = 161
$ws.Range("F16").Value = 14

# Sheet 全部类型 (All Types) - update 想去人数 column F
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 3013
$ws.Range("F3").Value = 3013
$ws.Range("F6").Value = 6434
$ws.Range("F7").Value = 2556
$ws.Range("F8").Value = 655
$ws.Range("F9").Value = 88
$ws.Range("F10").Value = 3147
$ws.Range("F11").Value = 370
$ws.Range("F15").Value = 7714
$ws.Range("F16").Value = 386
$ws.Range("F19").Value = 6
$ws.Range("F20").Value = 262
$ws.Range("F21").Value = 21
$ws.Range("F22").Value = 9488
$ws.Range("F24").Value = 271
$ws.Range("F27").Value = 132
$ws.Range("F28").Value = 36
$ws.Range("F31").Value = 2627
$ws.Range("F32").Value = 2054
$ws.Range("F33").Value = 18
$ws.Range("F34").Value = 1492
$ws.Range("F35").Value = 803
$ws.Range("F37").Value = 3974
$ws.Range("F38").Value = 221
$ws.Range("F39").Value = 617
$ws.Range("F41").Value = 1200
$ws.Range("F42").Value = 115
$ws.Range("F43").Value = 258
$ws.Range("F45").Value = 76
$ws.Range("F49").Value = 68
$ws.Range("F50").Value = 16
